# "Add applicants in Exam"
# Replace the three sample applicants (Vasudha/Anjali/Zubeda) with the new
# group of applicants (Aamena/Akram/Nazrin) for group_id 43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New applicant data
# ---------------------------------------------------------------------
$names     = @("Aamena", "Akram", "Nazrin")
$emails    = @("aamenas@gmail.com", "akram@gmail.com", "nazzz@gmail.com")
$passwords = @("aamena@123", "Akram@123", "Nazrin@123")
$phones    = @("9876543210", "9123456789", "9988776655")
$dobs      = @("2002-05-14", "2001-09-22", "2003-01-10")
$genders   = @("Male", "Female", "Male")
$addresses = @("Ahmedabad, Gujarat", "Vadodara, Gujarat", "Surat, Gujarat")
$groupId   = 43

# ---------------------------------------------------------------------
# Fill Phone / DOB / Gender / Address first (columns D:G), row by row --
# matches the order the workbook's shared-string table was built in.
# Phone and DOB must stay as literal TEXT (not auto-converted to a
# number / date serial), so force a text number format before writing
# and then drop back to the Normal style afterwards so no stray custom
# style sticks to the cell.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 3; $i++) {
    $row = $i + 2

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $phones[$i]
    $dCell.Style = "Normal"

    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $dobs[$i]
    $eCell.Style = "Normal"

    $ws.Cells.Item($row, 6).Value = $genders[$i]
    $ws.Cells.Item($row, 7).Value = $addresses[$i]
}

# ---------------------------------------------------------------------
# Full name / Email (columns A:B)
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 3; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $emails[$i]
}

# ---------------------------------------------------------------------
# Password (column C)
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 3; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $passwords[$i]
}

# ---------------------------------------------------------------------
# group_id (column H) - plain number
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 3; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $groupId
}

# ---------------------------------------------------------------------
# Re-point the mailto hyperlinks (B2:B4 -> email, C2:C4 -> password) at
# the new values. Delete the old hyperlinks first so Hyperlinks.Add
# replaces them in place (same rId order) instead of appending dupes.
# ---------------------------------------------------------------------
$ws.Range("B2:C4").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Cells.Item(2, 2), "mailto:" + $emails[0])
$ws.Hyperlinks.Add($ws.Cells.Item(3, 2), "mailto:" + $emails[1])
$ws.Hyperlinks.Add($ws.Cells.Item(4, 2), "mailto:" + $emails[2])
$ws.Hyperlinks.Add($ws.Cells.Item(2, 3), "mailto:" + $passwords[0])
$ws.Hyperlinks.Add($ws.Cells.Item(3, 3), "mailto:" + $passwords[1])
$ws.Hyperlinks.Add($ws.Cells.Item(4, 3), "mailto:" + $passwords[2])

# ---------------------------------------------------------------------
# Column A is now the widest column (names are short but let's mirror
# the author widening it); drop the old auto best-fit widths on B:E.
# ---------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 20.44140625

# ---------------------------------------------------------------------
# Selection moved to H9 in the saved file.
# ---------------------------------------------------------------------
$ws.Range("H9").Select()
